# Apply the IG-publisher "re-render" edits described by the diff:
#  - bump Version 5.0.0 -> 6.0.0
#  - bump Date to new publish timestamp
#  - replace the duplicated "Contact / No display for ContactDetail" rows
#    with "Publisher = Alvearie Team" and "Jurisdiction = United States of America"
#  - on the Elements sheet, the root Extension row's Short/Definition now
#    show the profile's own Title/Description instead of the generic
#    "Extension" / "An Extension" placeholders

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsElem = $wb.Worksheets.Item("Elements")

# --- Metadata sheet ---------------------------------------------------

# Version
$wsMeta.Range("B3").Value = "6.0.0"

# Date
$wsMeta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher (row 9) was blank -> now has a value
$wsMeta.Range("B9").Value = "Alvearie Team"

# Row 10 was a duplicate "Contact" / "No display for ContactDetail" row;
# turn it into the new "Jurisdiction" / "United States of America" row
$wsMeta.Range("A10").Value = "Jurisdiction"
$wsMeta.Range("B10").Value = "United States of America"

# Row 11 was the second duplicate "Contact" row - remove it entirely so
# everything below (Description, Purpose, Copyright, ...) shifts up by one
$wsMeta.Rows.Item(11).Delete()

# --- Elements sheet -----------------------------------------------------

# Root Extension row: Short / Definition now reflect the profile's own
# Title / Description rather than the generic Extension placeholders
$wsElem.Range("K2").Value = "Substance Abuse Ambulatory Coverage Indicator"
$wsElem.Range("L2").Value = "Indicator of Substance Abuse (chemical dependency) Ambulatory benefit coverage for the member. This finer granularity of MHSA benefit coverage may be used in HEDIS reporting."

# The "Short" column (K) is a few characters wider now that it holds the
# longer title text, so its best-fit width grows
$wsElem.Columns.Item(11).ColumnWidth = 46
